# Update NATMI ligand-receptor pair statistics after recomputation
# following Dr Hou's advice: ligand/receptor-expressing cell counts
# (columns E and K) changed from 1 to 3, which cascades into the
# dependent average/total expression and specificity columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 13.46482166666667
$ws.Range("H2").Value = 40.394465
$ws.Range("I2").Value = 0.4580736409596084
$ws.Range("J2").Value = 0.4580736409596083
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.03687233333333333
$ws.Range("N2").Value = 0.110617
$ws.Range("O2").Value = 0.02376577678298649
$ws.Range("P2").Value = 0.02376577678298649
$ws.Range("Q2").Value = 0.4964793927672222
$ws.Range("R2").Value = 4.468314534905
$ws.Range("S2").Value = 0.01088647590121595
$ws.Range("T2").Value = 0.01088647590121595
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 13.46482166666667
$ws.Range("H3").Value = 40.394465
$ws.Range("I3").Value = 0.4580736409596084
$ws.Range("J3").Value = 0.4580736409596083
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.291491666666667
$ws.Range("N3").Value = 3.874475
$ws.Range("O3").Value = 0.8324209479669633
$ws.Range("P3").Value = 0.8324209479669634
$ws.Range("Q3").Value = 17.38970497565278
$ws.Range("R3").Value = 156.507344780875
$ws.Range("S3").Value = 0.3813100944462756
$ws.Range("T3").Value = 0.3813100944462756
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 13.46482166666667
$ws.Range("H4").Value = 40.394465
$ws.Range("I4").Value = 0.4580736409596084
$ws.Range("J4").Value = 0.4580736409596083
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.2231246666666667
$ws.Range("N4").Value = 0.669374
$ws.Range("O4").Value = 0.1438132752500502
$ws.Range("P4").Value = 0.1438132752500502
$ws.Range("Q4").Value = 3.004333846101111
$ws.Range("R4").Value = 27.03900461491
$ws.Range("S4").Value = 0.06587707061211681
$ws.Range("T4").Value = 0.06587707061211681
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.78463766666667
$ws.Range("H5").Value = 38.353913
$ws.Range("I5").Value = 0.4349337606763218
$ws.Range("J5").Value = 0.4349337606763218
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.03687233333333333
$ws.Range("N5").Value = 0.110617
$ws.Range("O5").Value = 0.02376577678298649
$ws.Range("P5").Value = 0.02376577678298649
$ws.Range("Q5").Value = 0.4713994215912222
$ws.Range("R5").Value = 4.242594794321
$ws.Range("S5").Value = 0.01033653867161833
$ws.Range("T5").Value = 0.01033653867161833
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.78463766666667
$ws.Range("H6").Value = 38.353913
$ws.Range("I6").Value = 0.4349337606763218
$ws.Range("J6").Value = 0.4349337606763218
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.291491666666667
$ws.Range("N6").Value = 3.874475
$ws.Range("O6").Value = 0.8324209479669633
$ws.Range("P6").Value = 0.8324209479669634
$ws.Range("Q6").Value = 16.51125300785278
$ws.Range("R6").Value = 148.601277070675
$ws.Range("S6").Value = 0.3620479733650201
$ws.Range("T6").Value = 0.3620479733650201
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.78463766666667
$ws.Range("H7").Value = 38.353913
$ws.Range("I7").Value = 0.4349337606763218
$ws.Range("J7").Value = 0.4349337606763218
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.2231246666666667
$ws.Range("N7").Value = 0.669374
$ws.Range("O7").Value = 0.1438132752500502
$ws.Range("P7").Value = 0.1438132752500502
$ws.Range("Q7").Value = 2.852568017829111
$ws.Range("R7").Value = 25.673112160462
$ws.Range("S7").Value = 0.06254924863968332
$ws.Range("T7").Value = 0.06254924863968331
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.144988333333333
$ws.Range("H8").Value = 9.434965
$ws.Range("I8").Value = 0.1069925983640697
$ws.Range("J8").Value = 0.1069925983640697
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.03687233333333333
$ws.Range("N8").Value = 0.110617
$ws.Range("O8").Value = 0.02376577678298649
$ws.Range("P8").Value = 0.02376577678298649
$ws.Range("Q8").Value = 0.1159630581561111
$ws.Range("R8").Value = 1.043667523405
$ws.Range("S8").Value = 0.002542762210152207
$ws.Range("T8").Value = 0.002542762210152206
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.144988333333333
$ws.Range("H9").Value = 9.434965
$ws.Range("I9").Value = 0.1069925983640697
$ws.Range("J9").Value = 0.1069925983640697
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.291491666666667
$ws.Range("N9").Value = 3.874475
$ws.Range("O9").Value = 0.8324209479669633
$ws.Range("P9").Value = 0.8324209479669634
$ws.Range("Q9").Value = 4.061726224263889
$ws.Range("R9").Value = 36.555536018375
$ws.Range("S9").Value = 0.08906288015566749
$ws.Range("T9").Value = 0.0890628801556675
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.144988333333333
$ws.Range("H10").Value = 9.434965
$ws.Range("I10").Value = 0.1069925983640697
$ws.Range("J10").Value = 0.1069925983640697
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.2231246666666667
$ws.Range("N10").Value = 0.669374
$ws.Range("O10").Value = 0.1438132752500502
$ws.Range("P10").Value = 0.1438132752500502
$ws.Range("Q10").Value = 0.7017244735455556
$ws.Range("R10").Value = 6.315520261910001
$ws.Range("S10").Value = 0.01538695599825003
$ws.Range("T10").Value = 0.01538695599825003
